$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.148.51'
$ws.Range('E2').Value = '  -3.65%  '
$ws.Range('D3').Value = '2.454.07'
$ws.Range('E3').Value = '  -2.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.56'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.78'
$ws.Range('E6').Value = '  -8.05%  '
$ws.Range('E7').Value = '  -3.00%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.492'
$ws.Range('E9').Value = '  -5.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.08'
$ws.Range('E10').Value = '  -7.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0772'
$ws.Range('E11').Value = '  -4.00%  '
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.93'
$ws.Range('E13').Value = '  -5.50%  '
$ws.Range('D14').Value = '2.834.58'
$ws.Range('E14').Value = '  -2.83%  '
$ws.Range('D15').Value = '2.457.87'
$ws.Range('E15').Value = '  -3.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.56'
$ws.Range('E16').Value = '  -5.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.775'
$ws.Range('E17').Value = '  -4.15%  '
$ws.Range('D18').Value = '41.122.04'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.25'
$ws.Range('E19').Value = '  -6.84%  '
$ws.Range('D20').Value = '0.0₃0911'
$ws.Range('E20').Value = '  -4.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.11'
$ws.Range('E21').Value = '  -9.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.51'
$ws.Range('E22').Value = '  -3.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.92'
$ws.Range('E23').Value = '  -3.67%  '
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('E26').Value = '  -7.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.74'
$ws.Range('E27').Value = '  -6.83%  '
$ws.Range('E28').Value = '  -5.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.53'
$ws.Range('E29').Value = '  -5.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.34'
$ws.Range('E30').Value = '  -8.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '150.82'
$ws.Range('E31').Value = '  -4.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.43'
$ws.Range('E32').Value = '  -5.27%  '
$ws.Range('E33').Value = '  -5.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.55'
$ws.Range('E34').Value = '  -3.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0731'
$ws.Range('E35').Value = '  -6.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.96'
$ws.Range('E36').Value = '  -5.98%  '
$ws.Range('E37').Value = '  -6.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.59'
$ws.Range('E38').Value = '  -6.98%  '
$ws.Range('E39').Value = '  -3.75%  '
$ws.Range('E40').Value = '  -8.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.11'
$ws.Range('E41').Value = '  -1.83%  '
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.86'
$ws.Range('E43').Value = '  -8.94%  '
$ws.Range('D44').Value = '1.964.06'
$ws.Range('E44').Value = '  -2.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0281'
$ws.Range('E45').Value = '  -6.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.99'
$ws.Range('E46').Value = '  -9.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.47'
$ws.Range('E47').Value = '  -4.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '70.23'
$ws.Range('E48').Value = '  -2.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '95.71'
$ws.Range('E49').Value = '  -5.39%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.175'
$ws.Range('E50').Value = '  -7.73%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.33'
$ws.Range('E51').Value = '  -7.24%  '
